$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" - Right/Wrong counts
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total" - total marks and summary text
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "90 / 112"
